$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New task entries for rows 25-28 (timesheet continuation)
$ws.Range("A25").Value = "UC01 færdig"
$ws.Range("C25").Value = 43965
$ws.Range("D25").Value = 0.54166666666666663
$ws.Range("E25").Value = 0.58333333333333337

$ws.Range("A26").Value = "AD01 med Toke"
$ws.Range("C26").Value = 43965
$ws.Range("D26").Value = 0.58333333333333337
$ws.Range("E26").Value = 0.60416666666666663

$ws.Range("A27").Value = "DOM01 med MIK"
$ws.Range("C27").Value = 43965
$ws.Range("D27").Value = 0.60416666666666663
$ws.Range("E27").Value = 0.64583333333333337

$ws.Range("A28").Value = "SSD02"
$ws.Range("C28").Value = 43965
$ws.Range("D28").Value = 0.64583333333333337
$ws.Range("E28").Value = 0.67708333333333337

# Update sheet view: scrolled position + active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("H22").Select() | Out-Null
